# Update weekly Fruta/Hortaliza price rows for Vega Monumental Concepcion - Tuna
# Refresh Fecha/Calidad/Volumen/Precio(min,max,prom)/Unidad/Precio-por-Kg values
# to reflect the latest weekly snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row {
    param($Row, $Fecha, $Calidad, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $Unidad, $PrecioKg)

    $ws.Range("D$Row").Value = $Fecha
    $ws.Range("L$Row").Value = $Calidad
    $ws.Range("M$Row").Value = $Volumen
    $ws.Range("N$Row").Value = $PrecioMin
    $ws.Range("O$Row").Value = $PrecioMax
    $ws.Range("P$Row").Value = $PrecioProm
    $ws.Range("Q$Row").Value = $Unidad
    $ws.Range("S$Row").Value = $PrecioKg
}

Set-Row 2  45084 "Primera" 100 20000 21000 20500 "$/caja 18 kilos granel" 1139
Set-Row 3  44819 "Primera" 100 25000 26000 25500 "$/caja 18 kilos granel" 1417
Set-Row 4  45044 "Primera" 100 17000 18000 17500 "$/caja 18 kilos" 972
Set-Row 5  44316 "Primera" 50  20000 20000 20000 "$/caja 18 kilos" 1111
Set-Row 8  45030 "Primera" 100 15000 16000 15500 "$/caja 18 kilos granel" 861
Set-Row 9  45002 "Primera" 100 12000 13000 12500 "$/caja 18 kilos" 694
Set-Row 10 45014 "Primera" 50  13000 14000 13600 "$/caja 18 kilos" 756
Set-Row 11 45014 "Segunda" 20  10000 10000 10000 "$/caja 18 kilos" 556
Set-Row 12 44516 "Primera" 100 33000 34000 33500 "$/caja 18 kilos" 1861
Set-Row 13 44280 "Primera" 100 14000 15000 14500 "$/caja 18 kilos" 806
Set-Row 14 44280 "Segunda" 50  12000 12000 12000 "$/caja 18 kilos" 667
Set-Row 15 44687 "Primera" 100 18000 19000 18500 "$/caja 18 kilos" 1028

Write-Host "Rows updated"
